$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1400
$ws.Range("F3").Value = 210000

$ws.Range("D4").Value = 900
$ws.Range("F4").Value = 180000

$ws.Range("D13").Value = 80
$ws.Range("F13").Value = 32000

$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 'Asam mefenamat'
$ws.Range("E19").Value = 350
$ws.Range("F19").Value = 21000

$ws.Range("B20").Value = 10
$ws.Range("C20").Value = 'Charm wings'
$ws.Range("E20").Value = 900
$ws.Range("F20").Value = 54000

$ws.Range("B21").Value = 25
$ws.Range("C21").Value = 'Hansaplast '
$ws.Range("D21").Value = 60
$ws.Range("E21").Value = 400
$ws.Range("F21").Value = 24000

$ws.Range("B22").Value = 26
$ws.Range("C22").Value = 'Hansaplast koyo'
$ws.Range("D22").Value = 60
$ws.Range("E22").Value = 750
$ws.Range("F22").Value = 45000

$ws.Range("B23").Value = 82
$ws.Range("C23").Value = 'Ultrafix 5 cm x 1 m'
$ws.Range("D23").Value = 60
$ws.Range("E23").Value = 4200
$ws.Range("F23").Value = 252000

$ws.Range("B25").Value = 60
$ws.Range("C25").Value = 'Sarung tangan Uk M'
$ws.Range("E25").Value = 980
$ws.Range("F25").Value = 49000

$ws.Range("B26").Value = 86
$ws.Range("C26").Value = 'Verban 10 cm '
$ws.Range("D26").Value = 50
$ws.Range("E26").Value = 1100
$ws.Range("F26").Value = 55000

$ws.Range("B32").Value = 4
$ws.Range("C32").Value = 'Ambroxol '
$ws.Range("E32").Value = 300
$ws.Range("F32").Value = 9000

$ws.Range("B33").Value = 13
$ws.Range("C33").Value = 'Diapet'
$ws.Range("D33").Value = 30
$ws.Range("F33").Value = 9000

$ws.Range("B34").Value = 66
$ws.Range("C34").Value = 'Spuit 10 cc'
$ws.Range("D34").Value = 30
$ws.Range("E34").Value = 4000
$ws.Range("F34").Value = 120000

$ws.Range("B35").Value = 21
$ws.Range("C35").Value = 'Gelang tangan dewasa perempuan'
$ws.Range("D35").Value = 29
$ws.Range("E35").Value = 1000
$ws.Range("F35").Value = 29000

$ws.Range("B42").Value = 62
$ws.Range("C42").Value = 'Selang NGT '
$ws.Range("E42").Value = 29500
$ws.Range("F42").Value = 472000

$ws.Range("B43").Value = 68
$ws.Range("C43").Value = 'Spuit 20 '
$ws.Range("D43").Value = 16
$ws.Range("E43").Value = 8000
$ws.Range("F43").Value = 128000

$ws.Range("B44").Value = 58
$ws.Range("C44").Value = 'Sarung tangan panjang sterill'
$ws.Range("E44").Value = 20000
$ws.Range("F44").Value = 300000

$ws.Range("B47").Value = 5
$ws.Range("C47").Value = 'Aquades '

$ws.Range("D48").Value = 9
$ws.Range("F48").Value = 7200

$ws.Range("B49").Value = 46
$ws.Range("C49").Value = 'NaCl 500ml'
$ws.Range("D49").Value = 9
$ws.Range("E49").Value = 13000
$ws.Range("F49").Value = 117000

$ws.Range("B50").Value = 63
$ws.Range("C50").Value = 'Selang OGT'
$ws.Range("D50").Value = 7
$ws.Range("E50").Value = 17500
$ws.Range("F50").Value = 122500

$ws.Range("B51").Value = 57
$ws.Range("C51").Value = 'Sarung tangan obgyn sterill'
$ws.Range("E51").Value = 27000
$ws.Range("F51").Value = 162000

$ws.Range("B52").Value = 79
$ws.Range("C52").Value = 'Transfusi set '
$ws.Range("D52").Value = 6
$ws.Range("E52").Value = 11000
$ws.Range("F52").Value = 66000

$ws.Range("B53").Value = 7
$ws.Range("C53").Value = 'Assering 500 ml'
$ws.Range("E53").Value = 9500
$ws.Range("F53").Value = 47500

$ws.Range("B54").Value = 19
$ws.Range("C54").Value = 'Gelang tangan anak perempuan'
$ws.Range("E54").Value = 800
$ws.Range("F54").Value = 4000

$ws.Range("B55").Value = 55
$ws.Range("C55").Value = 'Sabun cuci Kimia '
$ws.Range("D55").Value = 5
$ws.Range("E55").Value = 15000
$ws.Range("F55").Value = 75000

$ws.Range("B56").Value = 16
$ws.Range("C56").Value = 'Extention tube'
$ws.Range("E56").Value = 12487
$ws.Range("F56").Value = 49948

$ws.Range("B57").Value = 40
$ws.Range("C57").Value = 'Masker nebulizer '
$ws.Range("D57").Value = 4
$ws.Range("E57").Value = 9000
$ws.Range("F57").Value = 36000

$ws.Range("B58").Value = 81
$ws.Range("C58").Value = 'Ultrafix 10 cm x 5 m'
$ws.Range("D58").Value = 4
$ws.Range("E58").Value = 32000
$ws.Range("F58").Value = 128000

$ws.Range("B59").Value = 29
$ws.Range("C59").Value = 'Inerasit gel'
$ws.Range("E59").Value = 125000
$ws.Range("F59").Value = 375000

$ws.Range("B60").Value = 33
$ws.Range("C60").Value = 'Insto'
$ws.Range("E60").Value = 16000
$ws.Range("F60").Value = 48000

$ws.Range("B61").Value = 36
$ws.Range("C61").Value = 'KA-EN 3A 500ml'
$ws.Range("D61").Value = 3
$ws.Range("E61").Value = 9500
$ws.Range("F61").Value = 28500

$ws.Range("B62").Value = 88
$ws.Range("C62").Value = 'Wool'
$ws.Range("D62").Value = 3
$ws.Range("E62").Value = 40000
$ws.Range("F62").Value = 120000

$ws.Range("B69").Value = 1
$ws.Range("C69").Value = 'Alginate'

$ws.Range("B70").Value = 3
$ws.Range("C70").Value = 'Allevyn '
$ws.Range("E70").Value = 125000
$ws.Range("F70").Value = 125000

$ws.Range("B71").Value = 8
$ws.Range("C71").Value = 'Betaplast'
$ws.Range("E71").Value = 40000
$ws.Range("F71").Value = 40000

$ws.Range("B72").Value = 14
$ws.Range("C72").Value = 'Excle care hydrosoloid'
$ws.Range("E72").Value = 50000
$ws.Range("F72").Value = 50000

$ws.Range("B73").Value = 15
$ws.Range("C73").Value = 'Exelcare'
$ws.Range("E73").Value = 130000
$ws.Range("F73").Value = 130000

$ws.Range("B74").Value = 28
$ws.Range("C74").Value = 'Indosorb'
$ws.Range("E74").Value = 100000
$ws.Range("F74").Value = 100000

$ws.Range("B75").Value = 32
$ws.Range("C75").Value = 'Inomed foam '
$ws.Range("E75").Value = 50000
$ws.Range("F75").Value = 50000

$ws.Range("B76").Value = 38
$ws.Range("C76").Value = 'Kilbac'
$ws.Range("E76").Value = 32000
$ws.Range("F76").Value = 32000

$ws.Range("B77").Value = 39
$ws.Range("C77").Value = 'Lauret '
$ws.Range("E77").Value = 175000
$ws.Range("F77").Value = 175000

$ws.Range("B78").Value = 41
$ws.Range("C78").Value = 'Medivix'
$ws.Range("E78").Value = 15000
$ws.Range("F78").Value = 15000

$ws.Range("B79").Value = 43
$ws.Range("C79").Value = 'Metcovazin green 10 gr '
$ws.Range("E79").Value = 40000
$ws.Range("F79").Value = 40000

$ws.Range("B80").Value = 44
$ws.Range("C80").Value = 'Metcovazin red 10 gr'
$ws.Range("E80").Value = 50000
$ws.Range("F80").Value = 50000

$ws.Range("B81").Value = 45
$ws.Range("C81").Value = 'Metcovazin red 25 gr'
$ws.Range("E81").Value = 80000
$ws.Range("F81").Value = 80000

$ws.Range("B82").Value = 50
$ws.Range("C82").Value = 'Ocpenic'
$ws.Range("E82").Value = 85000
$ws.Range("F82").Value = 85000

$ws.Range("B83").Value = 72
$ws.Range("C83").Value = 'Star AG'
$ws.Range("E83").Value = 250000
$ws.Range("F83").Value = 250000

$ws.Range("B84").Value = 74
$ws.Range("C84").Value = 'Suprasorb A'
$ws.Range("E84").Value = 70000
$ws.Range("F84").Value = 70000

$ws.Range("B85").Value = 75
$ws.Range("C85").Value = 'Tegaderm 10 x 12 '
$ws.Range("E85").Value = 25000
$ws.Range("F85").Value = 25000

$ws.Range("B86").Value = 76
$ws.Range("C86").Value = 'Tegaderm 6 x 7'
$ws.Range("D86").Value = 1
$ws.Range("E86").Value = 12000
$ws.Range("F86").Value = 12000

$ws.Range("B87").Value = 80
$ws.Range("C87").Value = 'Transparant film '
$ws.Range("D87").Value = 1
$ws.Range("E87").Value = 45000
$ws.Range("F87").Value = 45000

$ws.Range("B88").Value = 83
$ws.Range("C88").Value = 'Urgotul silver '
$ws.Range("D88").Value = 1
$ws.Range("E88").Value = 45000

$ws.Range("B89").Value = 87
$ws.Range("C89").Value = 'W care AH '
$ws.Range("D89").Value = 1
$ws.Range("E89").Value = 75000
